# Fix for the 20 minute trade bug: append a new trade row (row 4) so a
# trader no longer needs to enter the close price when data can't be
# located from Yahoo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 10015.91
$ws.Range("B4").Value = 10039
$ws.Range("C4").Value = 286.39
$ws.Range("D4").Value = 287.04000000000002
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = 0.23

$ws.Range("G4").Value = 42608.639652777776
$ws.Range("G4").NumberFormat = "m/d/yy h:mm"

$ws.Range("H4").Value = $false

# Column A's best-fit width shifts slightly (8.85546875 -> 9 chars) now that
# it holds an 8-character value (10015.91) in row 4.
$ws.Columns("A:A").ColumnWidth = 8.2
